# Apply the view-state + data edits captured by the diff.
#
# (Lower-level incidental metadata from the diff -- fileVersion/rupBuild,
# xr:revisionPtr GUIDs, workbookView x/y/window size, font-table locale
# names, theme name, sheetFormatPr default sizes and x14ac:dyDescent --
# are artifacts of the particular Excel build/session that last saved the
# file and are not reachable through the Excel object model exposed here,
# so they are intentionally left alone.)

$wb = $excel.ActiveWorkbook

# --- Sheet1 ("Sheet1"): new selection + zoom -------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate() | Out-Null
$wb.Windows.Item(1).Zoom = 193
$ws1.Range("D16").Select() | Out-Null

# --- Sheet3 ("12"): add a data point, then set selection + zoom ------------
# Doing this last keeps it the active/tabSelected sheet, matching
# workbook.xml's unchanged activeTab="2".
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(3, 4).Value = 1

$ws3.Activate() | Out-Null
$wb.Windows.Item(1).Zoom = 226
$ws3.Range("D6").Select() | Out-Null
